# Sample Project / Main.xlsx edit:
# Row 11, column B ("R40" rule-name label) is replaced with the literal
# text "1". It must remain a *text* value (not get auto-converted to a
# number) and must keep its existing cell style, so we build it as a
# formula that evaluates to the text "1" and then convert that formula
# result in place to a static value via Copy / PasteSpecial(Values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B11")
$cell.Formula = "=""1"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$excel.CutCopyMode = $false
